$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - "Save" (copy format from G1 so it reuses the same bold/border style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"
$excel.CutCopyMode = $false

# New column H data values (Save)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
